$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the edited cells keep Text format so numeric-looking strings
# (prices like "11.50", "1.00", "0.0000200") are not auto-converted to
# numbers on assignment, matching the source data which stores these as
# literal text (inline strings) rather than numeric values.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "96.440.95"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "3.575.56"
$ws.Range("E3").Value = "  -0.39%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "240.93"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "656.04"
$ws.Range("E6").Value = "  +0.94%  "

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "1.55"
$ws.Range("E7").Value = "  +6.08%  "

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.405"
$ws.Range("E8").Value = "  -1.10%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.83%  "

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "3.570.41"
$ws.Range("E11").Value = "  -0.51%  "

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "43.24"
$ws.Range("E12").Value = "  -0.56%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.93%  "

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "6.36"
$ws.Range("E14").Value = "  +0.27%  "

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "4.246.52"
$ws.Range("E15").Value = "  -0.69%  "

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "96.562.46"
$ws.Range("E16").Value = "  +1.21%  "

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000257"
$ws.Range("E17").Value = "  +0.03%  "

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "3.587.02"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "7.76"
$ws.Range("E19").Value = "  -2.48%  "

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "12.58"
$ws.Range("E20").Value = "  +1.19%  "

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "17.75"
$ws.Range("E21").Value = "  -1.86%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.01%  "

$ws.Range("B23:E23").NumberFormat = "@"
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "511.42"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("B24:E24").NumberFormat = "@"
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").Value = "3.44"
$ws.Range("E24").Value = "  -1.36%  "

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000200"
$ws.Range("E25").Value = "  +1.53%  "

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "6.81"
$ws.Range("E26").Value = "  +2.11%  "

$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "96.36"
$ws.Range("E27").Value = "  -0.23%  "

$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = "12.75"
$ws.Range("E28").Value = "  -1.70%  "

$ws.Range("D29:E29").NumberFormat = "@"
$ws.Range("D29").Value = "3.767.80"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("D30:E30").NumberFormat = "@"
$ws.Range("D30").Value = "3.00"
$ws.Range("E30").Value = "  -3.99%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.32%  "

$ws.Range("D32:E32").NumberFormat = "@"
$ws.Range("D32").Value = "11.50"
$ws.Range("E32").Value = "  +1.84%  "

$ws.Range("D33:E33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("D34:E34").NumberFormat = "@"
$ws.Range("D34").Value = "0.186"
$ws.Range("E34").Value = "  +5.47%  "

$ws.Range("D35:E35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.32%  "

$ws.Range("D36:E36").NumberFormat = "@"
$ws.Range("D36").Value = "31.60"
$ws.Range("E36").Value = "  -1.04%  "

$ws.Range("D37:E37").NumberFormat = "@"
$ws.Range("D37").Value = "0.563"
$ws.Range("E37").Value = "  +0.18%  "

$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = "599.01"
$ws.Range("E38").Value = "  +6.65%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.20%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +9.17%  "

$ws.Range("B41:E41").NumberFormat = "@"
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("B42:E42").NumberFormat = "@"
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.151"
$ws.Range("E42").Value = "  +0.29%  "

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.906"
$ws.Range("E43").Value = "  -2.27%  "

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "1.83"
$ws.Range("E44").Value = "  +5.53%  "

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "5.71"
$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "34.29"
$ws.Range("E46").Value = "  +1.81%  "

$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "2.27"
$ws.Range("E47").Value = "  +0.29%  "

$ws.Range("B48:E48").NumberFormat = "@"
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "23.52"
$ws.Range("E48").Value = "  -1.05%  "

$ws.Range("B49:E49").NumberFormat = "@"
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0419"
$ws.Range("E49").Value = "  -0.14%  "

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "3.64"
$ws.Range("E50").Value = "  +5.50%  "

$ws.Range("B51:E51").NumberFormat = "@"
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "53.51"
$ws.Range("E51").Value = "  -1.52%  "
